# review parte fabio 2
# Applies three textual revisions to the document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1 (SLIDE 8 paragraph): expand the Google Cardboard description
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(22)

$r = $d.Range($p1.Range.Start, $p1.Range.End)
$r.Find.ClearFormatting()
$r.Find.Replacement.ClearFormatting()
$ok1 = $r.Find.Execute("economica. Questo consente", $false, $false, $false, $false, $false, $true, 1, $false, "economica. Essa è formata da un cartone e due lenti e consente la visione stereoscopica. La sua economicità consente", 2)
Write-Host "change1 step1:" $ok1

$p1 = $d.Paragraphs.Item(22)
$r = $d.Range($p1.Range.Start, $p1.Range.End)
$r.Find.ClearFormatting()
$r.Find.Replacement.ClearFormatting()
$ok2 = $r.Find.Execute("collegato ed hanno prezzi", $false, $false, $false, $false, $false, $true, 1, $false, "collegato e sono caratterizzati da prezzi", 2)
Write-Host "change1 step2:" $ok2

$p1 = $d.Paragraphs.Item(22)
$r = $d.Range($p1.Range.Start, $p1.Range.End)
$r.Find.ClearFormatting()
$r.Find.Replacement.ClearFormatting()
$r.Find.Replacement.Font.Bold = 1
$ok3 = $r.Find.Execute("trattamento", $true, $false, $false, $false, $false, $true, 1, $false, "trattamento", 2)
Write-Host "change1 step3 (bold trattamento):" $ok3

# ---------------------------------------------------------------------
# Change 2 (SLIDE 9 paragraph): rewrite the closing sentence
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(25)

$r = $d.Range($p2.Range.Start, $p2.Range.End)
$r.Find.ClearFormatting()
$r.Find.Replacement.ClearFormatting()
$ok4 = $r.Find.Execute("in modo innovativo e coinvolgente, sfruttando inoltre intrisecamente l’intattenimento volontario dato dal gioco.", $false, $false, $false, $false, $false, $true, 1, $false, "in modo interattivo e coinvolgente. Inoltre, l’intrattenimento volontario dato dal gioco consente un metodo di cura innovativo.", 2)
Write-Host "change2 step1:" $ok4

$p2 = $d.Paragraphs.Item(25)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$r2.Find.ClearFormatting()
$ok5 = $r2.Find.Execute("coinvolgente. Inoltre", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "change2 step2 locate:" $ok5
$boldR = $d.Range($r2.Start + 12, $r2.Start + 14)
Write-Host "change2 step2 bold text [" $boldR.Text "]"
$boldR.Font.Bold = 1

$p2 = $d.Paragraphs.Item(25)
$r3 = $d.Range($p2.Range.Start, $p2.Range.End)
$r3.Find.ClearFormatting()
$ok6 = $r3.Find.Execute("Inoltre,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "change2 step3 locate:" $ok6
$commaR = $d.Range($r3.End - 1, $r3.End)
Write-Host "change2 step3 bold text [" $commaR.Text "]"
$commaR.Font.Bold = 1

$p2 = $d.Paragraphs.Item(25)
$r4 = $d.Range($p2.Range.Start, $p2.Range.End)
$r4.Find.ClearFormatting()
$r4.Find.Replacement.ClearFormatting()
$r4.Find.Replacement.Font.Bold = 1
$ok7 = $r4.Find.Execute("innovativo", $true, $false, $false, $false, $false, $true, 1, $false, "innovativo", 2)
Write-Host "change2 step4 (bold innovativo):" $ok7

# ---------------------------------------------------------------------
# Change 3 (SLIDE 13 paragraph): rewrite the GameThread sentence and
# move the _GoBack bookmark so it still sits right before "a cui si
# accede attraverso il metodo runOnUIThread"
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(28)
$r = $d.Range($p3.Range.Start, $p3.Range.End)
$r.Find.ClearFormatting()
$r.Find.Replacement.ClearFormatting()
$ok8 = $r.Find.Execute("essa vengono effettuate le principali operazioni di aggiornamento del gioco, come la gestione delle animazioni delle auto nemiche e l’aggiornamento dei punteggi. Questo è stato fatto per non gravare sul Thread principale dell’applicazione, a cui si accede attraverso il metodo runOnUIThread. ", $false, $false, $false, $false, $false, $true, 1, $false, "essa vengono effettuate le principali operazioni di aggiornamento del gioco, come la gestione delle animazioni delle auto nemiche e l’aggiornamento dei punteggi. Ciò è stato fatto per non gravare  e sovraccaricare il Thread principale dell’applicazione, a cui si accede attraverso il metodo runOnUIThread. ", 2)
Write-Host "change3 step1:" $ok8

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$p3 = $d.Paragraphs.Item(28)
$r2 = $d.Range($p3.Range.Start, $p3.Range.End)
$r2.Find.ClearFormatting()
$ok9 = $r2.Find.Execute("dell’applicazione,", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "change3 step2 locate bookmark anchor:" $ok9
$d.Bookmarks.Add("_GoBack", $d.Range($r2.End, $r2.End))
Write-Host "bookmark re-added"
